$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -13.101
$ws.Range("E3").Value = 16.228
$ws.Range("A12").Value = -21.651
$ws.Range("C14").Value = -12.586
$ws.Range("E20").Value = 16.306
$ws.Range("E25").Value = 16.942
$ws.Range("C26").Value = -12.393
$ws.Range("A27").Value = -21.876
$ws.Range("E30").Value = 16.108
$ws.Range("C31").Value = -13.306
$ws.Range("A32").Value = -21.695
$ws.Range("C35").Value = -12.762
$ws.Range("A36").Value = -20.178
$ws.Range("C37").Value = -13.329
$ws.Range("A38").Value = -19.741
$ws.Range("E44").Value = 16.611
$ws.Range("C45").Value = -12.702
$ws.Range("A46").Value = -21.879
$ws.Range("E47").Value = 16.209
$ws.Range("C52").Value = -11.286
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.145
$ws.Range("A56").Value = -22.061
$ws.Range("C57").Value = -13.829
$ws.Range("E58").Value = 16.593
$ws.Range("A67").Value = -21.554
$ws.Range("A69").Value = -21.656
$ws.Range("A72").Value = -21.431
$ws.Range("E78").Value = 16.482
$ws.Range("C81").Value = -13.427
$ws.Range("A83").Value = -21.627
$ws.Range("C83").Value = -13.026
$ws.Range("E84").Value = 16.399
$ws.Range("A86").Value = -22.257
$ws.Range("E89").Value = 17.108
$ws.Range("A91").Value = -21.508
$ws.Range("E91").Value = 17.39
$ws.Range("E92").Value = 17.293
$ws.Range("A93").Value = -21.421
$ws.Range("E96").Value = 16.265
$ws.Range("A99").Value = -20.428
$ws.Range("C100").Value = -12.578
$ws.Range("C102").Value = -12.87
$ws.Range("E102").Value = 16.418
